# Update with new TPM values (re-run of NATMI lrc2p scoring for Inhba-Acvr1b).
#
# The worksheet lists, for every (sending cluster -> target cluster) pair, the
# ligand (Inhba) stats for the sending cluster and the receptor (Acvr1b) stats
# for the target cluster, plus edge-level scores derived from them:
#   F  = Ligand detection rate               = E / TotalCells
#   H  = Ligand total expression value       = G * TotalCells
#   I  = Ligand avg-expr specificity         = G / sum(G over sending clusters)
#   J  = Ligand total-expr specificity       = H / sum(H over sending clusters)
#   L  = Receptor detection rate             = K / TotalCells        (unchanged)
#   N  = Receptor total expression value     = M * TotalCells
#   O  = Receptor avg-expr specificity       = M / sum(M over target clusters)
#   P  = Receptor total-expr specificity     = N / sum(N over target clusters)
#   Q  = Edge avg expression weight          = G * M
#   R  = Edge total expression weight        = H * N
#   S  = Edge avg expression specificity     = I * O
#   T  = Edge total expression specificity   = J * P
#
# The new TPM recomputation changed the underlying per-cluster ligand values
# (Ligand-expressing cells E, Ligand average expression G) for the ECs and
# Resolving-Mac sending clusters, and the per-cluster receptor average
# expression M for the ECs, MuSCs and Resolving-Mac target clusters. Every
# column derived from those (F,H,I,J,N,O,P,Q,R,S,T) is recalculated below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$TotalCells = 3

# Sending clusters, in the order they appear in the sheet (row blocks 2-5, 6-9, 10-13, 14-17)
$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# New ligand-expressing cell counts and ligand average expression values, per sending cluster
$E = @{ "ECs" = 2; "FAPs" = 3; "MuSCs" = 3; "Resolving-Mac" = 1 }
$G = @{ "ECs" = 0.1285113333333333; "FAPs" = 2.588894; "MuSCs" = 0.692415; "Resolving-Mac" = 0.01787866666666667 }

# New receptor average expression values, per target cluster
$M = @{ "ECs" = 3.483060666666667; "FAPs" = 4.620706999999999; "MuSCs" = 2.773309666666667; "Resolving-Mac" = 2.902635666666666 }

# Derived per-cluster totals
$H = @{}
$N = @{}
foreach ($c in $clusters) {
    $H[$c] = $G[$c] * $TotalCells
    $N[$c] = $M[$c] * $TotalCells
}

$sumG = 0; $sumH = 0; $sumM = 0; $sumN = 0
foreach ($c in $clusters) {
    $sumG += $G[$c]
    $sumH += $H[$c]
    $sumM += $M[$c]
    $sumN += $N[$c]
}

$I = @{}
$J = @{}
$O = @{}
$P = @{}
foreach ($c in $clusters) {
    $I[$c] = $G[$c] / $sumG
    $J[$c] = $H[$c] / $sumH
    $O[$c] = $M[$c] / $sumM
    $P[$c] = $N[$c] / $sumN
}

# Walk the 4x4 block of rows (sending cluster x target cluster), rows 2-17
$row = 2
foreach ($a in $clusters) {
    foreach ($d in $clusters) {
        $ws.Cells.Item($row, 5).Value2  = $E[$a]                     # E: Ligand-expressing cells
        $ws.Cells.Item($row, 6).Value2  = $E[$a] / $TotalCells        # F: Ligand detection rate
        $ws.Cells.Item($row, 7).Value2  = $G[$a]                      # G: Ligand average expression value
        $ws.Cells.Item($row, 8).Value2  = $H[$a]                      # H: Ligand total expression value
        $ws.Cells.Item($row, 9).Value2  = $I[$a]                      # I: Ligand derived specificity of average
        $ws.Cells.Item($row, 10).Value2 = $J[$a]                      # J: Ligand derived specificity of total

        $ws.Cells.Item($row, 13).Value2 = $M[$d]                      # M: Receptor average expression value
        $ws.Cells.Item($row, 14).Value2 = $N[$d]                      # N: Receptor total expression value
        $ws.Cells.Item($row, 15).Value2 = $O[$d]                      # O: Receptor derived specificity of average
        $ws.Cells.Item($row, 16).Value2 = $P[$d]                      # P: Receptor derived specificity of total

        $ws.Cells.Item($row, 17).Value2 = $G[$a] * $M[$d]              # Q: Edge average expression weight
        $ws.Cells.Item($row, 18).Value2 = $H[$a] * $N[$d]              # R: Edge total expression weight
        $ws.Cells.Item($row, 19).Value2 = $I[$a] * $O[$d]              # S: Edge average expression derived specificity
        $ws.Cells.Item($row, 20).Value2 = $J[$a] * $P[$d]              # T: Edge total expression derived specificity

        $row++
    }
}
